$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 0.5597564161496535
$ws.Range("P2").Value = 0.5597564161496534
$ws.Range("S2").Value = 0.5597564161496535
$ws.Range("T2").Value = 0.5597564161496534

# Row 3 updates
$ws.Range("M3").Value = 0.4111863333333334
$ws.Range("N3").Value = 1.233559
$ws.Range("O3").Value = 0.4402435838503465
$ws.Range("P3").Value = 0.4402435838503465
$ws.Range("Q3").Value = 0.01608574642211111
$ws.Range("R3").Value = 0.144771717799
$ws.Range("S3").Value = 0.4402435838503465
$ws.Range("T3").Value = 0.4402435838503465
